$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "TIPO DOCUMENTO"
$ws.Range("B1").Value = "NUMERO DOCUMENTO"
$ws.Range("C1").Value = "SUELDO"

# --- Data rows: filtered/cleaned export, document numbers as numbers ---
$ws.Range("A2").Value = "CC"
$ws.Range("B2").Value = 1076669247
$ws.Range("C2").Value = 1500000

$ws.Range("A3").Value = "CC"
$ws.Range("B3").Value = 1005043696
$ws.Range("C3").Value = 900000

$ws.Range("A4").Value = "CEDULA"
$ws.Range("B4").Value = 74320814
$ws.Range("C4").Value = 1200000

$ws.Range("A5").Value = "CE"
$ws.Range("B5").Value = 39743595
$ws.Range("C5").Value = -500000

$ws.Range("A6").Value = "CE"
$ws.Range("B6").Value = 72257456
$ws.Range("C6").Value = 344324234

# --- New row appended from the export ---
$ws.Range("A7").Value = "CC"
$ws.Range("B7").Value = 72257456
$ws.Range("C7").Value = 458000

# Widen the SUELDO column like the resulting export
$ws.Columns.Item(3).ColumnWidth = 11.43

# Leave selection where the export left off
$ws.Range("D10").Select()
